# Update "想去人数" (want-to-go count) values and one "是否有舞台" flag
# across the 展览 (sheet1), 本地生活 (sheet3) and 全部类型 (sheet4) worksheets,
# reflecting refreshed data generated at commit a3196b5.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 971
$ws1.Range("F6").Value  = 6733
$ws1.Range("F8").Value  = 906
$ws1.Range("F10").Value = 746
$ws1.Range("F11").Value = 509
$ws1.Range("F13").Value = 40
$ws1.Range("F14").Value = 366
$ws1.Range("F15").Value = 816
$ws1.Range("F16").Value = 2413
$ws1.Range("F17").Value = 90
$ws1.Range("F18").Value = 170
$ws1.Range("F19").Value = 733
$ws1.Range("F20").Value = 30
$ws1.Range("F21").Value = 398
$ws1.Range("F22").Value = 39
$ws1.Range("F23").Value = 176
$ws1.Range("F24").Value = 14
$ws1.Range("F25").Value = 85
$ws1.Range("F26").Value = 19
$ws1.Range("H26").Value = $true
$ws1.Range("F27").Value = 116
$ws1.Range("F28").Value = 8
$ws1.Range("F29").Value = 9
$ws1.Range("F30").Value = 289
$ws1.Range("F31").Value = 7
$ws1.Range("F32").Value = 225

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 181

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 181
$ws4.Range("F3").Value  = 971
$ws4.Range("F10").Value = 6733
$ws4.Range("F12").Value = 906
$ws4.Range("F14").Value = 746
$ws4.Range("F15").Value = 509
$ws4.Range("F17").Value = 40
$ws4.Range("F18").Value = 366
$ws4.Range("F19").Value = 816
$ws4.Range("F21").Value = 2413
$ws4.Range("F22").Value = 90
$ws4.Range("F24").Value = 170
$ws4.Range("F25").Value = 733
$ws4.Range("F26").Value = 30
$ws4.Range("F27").Value = 398
$ws4.Range("F28").Value = 39
$ws4.Range("F29").Value = 176
$ws4.Range("F30").Value = 14
$ws4.Range("F31").Value = 85
$ws4.Range("F32").Value = 19
$ws4.Range("H32").Value = $true
$ws4.Range("F33").Value = 116
$ws4.Range("F34").Value = 8
$ws4.Range("F35").Value = 9
$ws4.Range("F36").Value = 289
$ws4.Range("F37").Value = 7
$ws4.Range("F38").Value = 225
